# sat-ecf-local-sep-5-2017.xlsx — re-checked radiant-azimuth sign errors
# (AS, BC, Kosice, CB, SC) and flipped the sign on the AS / SC inputs so the
# ECF->local conversion matches the ground-based measurements better, then
# reran the sheet with the corrected ECF velocity components (X/Y/Z) and
# longitude/latitude inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Longitude= (B3) / Latitude= (B4)
$ws.Range("B3").Value = 40.200000000000003
$ws.Range("B4").Value = 39.1

# corrected Unc row value (B5)
$ws.Range("B5").Value = 39.799999999999997

# ECF velocity components X (F2) / Y (F3) / Z (F4)
$ws.Range("F2").Value = 10.3
$ws.Range("F3").Value = 12.2
$ws.Range("F4").Value = -18

# Latitude / Longitude reference values (F5) and Unc (F6)
$ws.Range("F5").Value = 39.1
$ws.Range("F6").Value = 40.200000000000003

# the downstream Vn/Vd/Ve, Radiant_azimuth, Zenith Distance and Speed cells
# (C6:C8, D10:D12) are formulas and recalculate automatically.

# match the author's final cursor position/selection
$ws.Activate()
$ws.Range("L21").Select()
